$wb = $excel.ActiveWorkbook

# OFF sheet - update Week 15/16 "R" row (row 3) target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 284
$wsOff.Range("C3").Value = 206
$wsOff.Range("D3").Value = 66
$wsOff.Range("E3").Value = 33
$wsOff.Range("F3").Value = 6

# DEF sheet - update Week 15/16 "R" row (row 3) target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 275
$wsDef.Range("C3").Value = 194
$wsDef.Range("D3").Value = 53
$wsDef.Range("E3").Value = 28
